$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.70123487346835
$ws.Range("C2").Value = 7.658386050022585
$ws.Range("E2").Value = 21.36261348978845
$ws.Range("F2").Value = 38.37922495233495
$ws.Range("G2").Value = 3.614849352230174
$ws.Range("J2").Value = 7.674993271905417
$ws.Range("M2").Value = 18.84418612527423
$ws.Range("N2").Value = 17.12616918317195
$ws.Range("O2").Value = 20.12388729016065
$ws.Range("B3").Value = 11.19032637952035
$ws.Range("C3").Value = 7.310391844228143
$ws.Range("E3").Value = 21.33120191511144
$ws.Range("F3").Value = 38.29901435393634
$ws.Range("G3").Value = 3.6168905678764
$ws.Range("J3").Value = 7.692949633881477
$ws.Range("M3").Value = 18.66590936621242
$ws.Range("N3").Value = 17.1804997178728
$ws.Range("O3").Value = 20.17093952112802
$ws.Range("B4").Value = 10.86579563516508
$ws.Range("C4").Value = 7.086697488044908
$ws.Range("E4").Value = 21.31578405168231
$ws.Range("F4").Value = 38.26021885541592
$ws.Range("G4").Value = 3.618210138490815
$ws.Range("J4").Value = 7.704521281693746
$ws.Range("M4").Value = 18.55880693066664
$ws.Range("N4").Value = 17.21570734213931
$ws.Range("O4").Value = 20.20557592714331
$ws.Range("B5").Value = 10.73101199899181
$ws.Range("C5").Value = 6.993094136585347
$ws.Range("E5").Value = 21.3104782690251
$ws.Range("F5").Value = 38.24704787699326
$ws.Range("G5").Value = 3.618764588419547
$ws.Range("J5").Value = 7.709374634003067
$ws.Range("M5").Value = 18.51579689897248
$ws.Range("N5").Value = 17.23052055912065
$ws.Range("O5").Value = 20.22112984580378
$ws.Range("B6").Value = 10.70848420410102
$ws.Range("C6").Value = 6.977406055803197
$ws.Range("E6").Value = 21.30965637410863
$ws.Range("F6").Value = 38.24502043460858
$ws.Range("G6").Value = 3.618857665480544
$ws.Range("J6").Value = 7.710188865389432
$ws.Range("M6").Value = 18.50869466543131
$ws.Range("N6").Value = 17.23300844729307
$ws.Range("O6").Value = 20.22379931918242
$ws.Range("B7").Value = 10.86398789084838
$ws.Range("C7").Value = 7.085444917053799
$ws.Range("E7").Value = 21.31570853463634
$ws.Range("F7").Value = 38.26003053307871
$ws.Range("G7").Value = 3.618217548250055
$ws.Range("J7").Value = 7.704586177146069
$ws.Range("M7").Value = 18.55822425687092
$ws.Range("N7").Value = 17.21590523091219
$ws.Range("O7").Value = 20.20577987285118
$ws.Range("B8").Value = 11.52742715376401
$ws.Range("C8").Value = 7.540519375186967
$ws.Range("E8").Value = 21.35098251240834
$ws.Range("F8").Value = 38.34940524777043
$ws.Range("G8").Value = 3.615539444276956
$ws.Range("J8").Value = 7.681071482179595
$ws.Range("M8").Value = 18.7822517227443
$ws.Range("N8").Value = 17.14451931319015
$ws.Range("O8").Value = 20.13891505095276
$ws.Range("B9").Value = 12.7354181834289
$ws.Range("C9").Value = 8.350672494196578
$ws.Range("E9").Value = 21.45064010215403
$ws.Range("F9").Value = 38.60706607230819
$ws.Range("G9").Value = 3.610810954836881
$ws.Range("J9").Value = 7.639274934073735
$ws.Range("M9").Value = 19.23842076996198
$ws.Range("N9").Value = 17.01915243793947
$ws.Range("O9").Value = 20.05361704171103
$ws.Range("B10").Value = 13.55828549208625
$ws.Range("C10").Value = 8.892818552334926
$ws.Range("E10").Value = 21.54211862825781
$ws.Range("F10").Value = 38.8456834192248
$ws.Range("G10").Value = 3.607652502080238
$ws.Range("J10").Value = 7.611170407558106
$ws.Range("M10").Value = 19.58127864298968
$ws.Range("N10").Value = 16.9358946931427
$ws.Range("O10").Value = 20.01916832114583
$ws.Range("B11").Value = 13.9172981980687
$ws.Range("C11").Value = 9.12747636933093
$ws.Range("E11").Value = 21.58761651063489
$ws.Range("F11").Value = 38.96469204896042
$ws.Range("G11").Value = 3.606283433363782
$ws.Range("J11").Value = 7.598944348209256
$ws.Range("M11").Value = 19.73837449443594
$ws.Range("N11").Value = 16.89992688180938
$ws.Range("O11").Value = 20.00967382478912
$ws.Range("B12").Value = 14.05095846604176
$ws.Range("C12").Value = 9.21458506304211
$ws.Range("E12").Value = 21.60539569132542
$ws.Range("F12").Value = 39.01123616226732
$ws.Range("G12").Value = 3.605774686451464
$ws.Range("J12").Value = 7.594394575822538
$ws.Range("M12").Value = 19.79797831814894
$ws.Range("N12").Value = 16.88658001839172
$ws.Range("O12").Value = 20.00696972160798
$ws.Range("B13").Value = 14.02227544046155
$ws.Range("C13").Value = 9.195902951917672
$ws.Range("E13").Value = 21.60154230011153
$ws.Range("F13").Value = 39.00114675642963
$ws.Range("G13").Value = 3.605883824032215
$ws.Range("J13").Value = 7.595370900134431
$ws.Range("M13").Value = 19.7851371997366
$ws.Range("N13").Value = 16.88944235858008
$ws.Range("O13").Value = 20.00751241716006
$ws.Range("B14").Value = 13.92834085940943
$ws.Range("C14").Value = 9.134678102697153
$ws.Range("E14").Value = 21.58906823956525
$ws.Range("F14").Value = 38.96849179557999
$ws.Range("G14").Value = 3.606241384534983
$ws.Range("J14").Value = 7.598568435173875
$ws.Range("M14").Value = 19.74327610050578
$ws.Range("N14").Value = 16.89882335311973
$ws.Range("O14").Value = 20.00943347899677
$ws.Range("B15").Value = 13.87050261089496
$ws.Range("C15").Value = 9.096947263811673
$ws.Range("E15").Value = 21.58149889961065
$ws.Range("F15").Value = 38.94868136820757
$ws.Range("G15").Value = 3.60646166109136
$ws.Range("J15").Value = 7.600537420815582
$ws.Range("M15").Value = 19.71764855453249
$ws.Range("N15").Value = 16.90460506034098
$ws.Range("O15").Value = 20.01072633387363
$ws.Range("B16").Value = 13.53450717439349
$ws.Range("C16").Value = 8.877239584332418
$ws.Range("E16").Value = 21.5392226867689
$ws.Range("F16").Value = 38.83811424563466
$ws.Range("G16").Value = 3.607743332555491
$ws.Range("J16").Value = 7.611980621141422
$ws.Range("M16").Value = 19.57103071966232
$ws.Range("N16").Value = 16.93828357050304
$ws.Range("O16").Value = 20.01991337297229
$ws.Range("B17").Value = 13.32439551031622
$ws.Range("C17").Value = 8.739367583070342
$ws.Range("E17").Value = 21.51427659663607
$ws.Range("F17").Value = 38.77294632742985
$ws.Range("G17").Value = 3.608546907797488
$ws.Range("J17").Value = 7.619143506735765
$ws.Range("M17").Value = 19.48134070294982
$ws.Range("N17").Value = 16.95943198710449
$ws.Range("O17").Value = 20.02713365168679
$ws.Range("B18").Value = 13.20210772693667
$ws.Range("C18").Value = 8.658942263979149
$ws.Range("E18").Value = 21.50029437410247
$ws.Range("F18").Value = 38.73644977303153
$ws.Range("G18").Value = 3.609015480764863
$ws.Range("J18").Value = 7.623316031892409
$ws.Range("M18").Value = 19.42986237107382
$ws.Range("N18").Value = 16.97177549675724
$ws.Range("O18").Value = 20.03186777099273
$ws.Range("B19").Value = 13.16045927235319
$ws.Range("C19").Value = 8.631519431476258
$ws.Range("E19").Value = 21.49562335768828
$ws.Range("F19").Value = 38.72426281852145
$ws.Range("G19").Value = 3.609175228415823
$ws.Range("J19").Value = 7.624737827115045
$ws.Range("M19").Value = 19.41245285433179
$ws.Range("N19").Value = 16.9759856509139
$ws.Range("O19").Value = 20.03357038046964
$ws.Range("B20").Value = 13.346911672874
$ws.Range("C20").Value = 8.75416096253743
$ws.Range("E20").Value = 21.51689431819719
$ws.Range("F20").Value = 38.77978166500916
$ws.Range("G20").Value = 3.608460706128271
$ws.Range("J20").Value = 7.618375561890252
$ws.Range("M20").Value = 19.49087742063523
$ws.Range("N20").Value = 16.95716213031763
$ws.Range("O20").Value = 20.02630486357977
$ws.Range("B21").Value = 13.9559945065547
$ws.Range("C21").Value = 9.15270906362667
$ws.Range("E21").Value = 21.59271731335865
$ws.Range("F21").Value = 38.97804345344606
$ws.Range("G21").Value = 3.606136097764123
$ws.Range("J21").Value = 7.597627074314667
$ws.Range("M21").Value = 19.75556897245802
$ws.Range("N21").Value = 16.89606051475098
$ws.Range("O21").Value = 20.00884500657683
$ws.Range("B22").Value = 14.34068578143682
$ws.Range("C22").Value = 9.402966115990013
$ws.Range("E22").Value = 21.64547420799788
$ws.Range("F22").Value = 39.11622080541496
$ws.Range("G22").Value = 3.604673288133353
$ws.Range("J22").Value = 7.584532691963355
$ws.Range("M22").Value = 19.92920910641255
$ws.Range("N22").Value = 16.8577201129819
$ws.Range("O22").Value = 20.00262989136314
$ws.Range("B23").Value = 14.13661859747913
$ws.Range("C23").Value = 9.270342677172509
$ws.Range("E23").Value = 21.61702679273261
$ws.Range("F23").Value = 39.04169514448498
$ws.Range("G23").Value = 3.605448867632773
$ws.Range("J23").Value = 7.591478903020087
$ws.Range("M23").Value = 19.83649021543875
$ws.Range("N23").Value = 16.878037599995
$ws.Range("O23").Value = 20.00547073050395
$ws.Range("B24").Value = 13.33673675554251
$ws.Range("C24").Value = 8.747476489232058
$ws.Range("E24").Value = 21.51570972545775
$ws.Range("F24").Value = 38.77668838742878
$ws.Range("G24").Value = 3.608499657381521
$ws.Range("J24").Value = 7.6187225799973
$ws.Range("M24").Value = 19.48656560014364
$ws.Range("N24").Value = 16.95818775621459
$ws.Range("O24").Value = 20.02667774261104
$ws.Range("B25").Value = 12.41947995447086
$ws.Range("C25").Value = 8.140646495900011
$ws.Range("E25").Value = 21.42044718787906
$ws.Range("F25").Value = 38.5286245261612
$ws.Range("G25").Value = 3.612034473159389
$ws.Range("J25").Value = 7.650122842593207
$ws.Range("M25").Value = 19.11347338886675
$ws.Range("N25").Value = 17.05150894703908
$ws.Range("O25").Value = 20.07175514454467
